$wb = $excel.ActiveWorkbook

# --- Sheet: Pansies Alive --- (W15 Thursday garden data collection)
$ws1 = $wb.Worksheets.Item("Pansies Alive")
$ws1.Range("B26").Value = 5
$ws1.Range("C26").Value = 7
$ws1.Range("D26").Value = 4
$ws1.Range("E26").Value = 5
$ws1.Range("F26").Value = 1
[void]$ws1.Range("H30").Select()

# --- Sheet: Pansies Dead ---
$ws2 = $wb.Worksheets.Item("Pansies Dead")
$ws2.Range("B26").Value = 0
$ws2.Range("C26").Value = 0
$ws2.Range("D26").Value = 0
$ws2.Range("E26").Value = 0
$ws2.Range("F26").Value = 0
[void]$ws2.Range("D27").Select()

# --- Sheet: Cardoon (1) ---
$ws3 = $wb.Worksheets.Item("Cardoon (1)")
$ws3.Range("G25").Value = 17
$ws3.Range("C26").Value = 19.75
$ws3.Range("F26").Value = 12
$ws3.Range("G26").Value = 17
$ws3.Range("H26").Value = 7.75
[void]$ws3.Range("G23").Select()

# --- Sheet: Cardoon (2) ---
$ws4 = $wb.Worksheets.Item("Cardoon (2)")
$ws4.Range("D26").Value = 27
$ws4.Range("E26").Value = 24.5
$ws4.Range("F26").Value = 24.5
$ws4.Range("G26").Value = 24.5
$ws4.Range("H26").Value = 15
$ws4.Activate()
[void]$ws4.Range("J24").Select()
